# "Generate Report for Handoff" - a35e0bd0 and 4d8f17aa swap places (row 6 <-> row 7)
# on every sheet, and a35e0bd0's "Latest HO Xliff Generate Date" / handoff timestamp is
# refreshed to a newer value since it is now the most-recently-handed-off file.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview": columns A..G
#   A File Name, B Path And Name, C Extension, D Publish URL,
#   E zh-cn, F de-de, G Latest HO Xliff Generate Date
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A6").Value = "4d8f17aa-978e-4fa6-a9bf-edb8ba655f6c.md"
$wsOverview.Range("B6").Value = "e2e\4d8f17aa-978e-4fa6-a9bf-edb8ba655f6c.md"
$wsOverview.Range("E6").Value = "Ready for handoff"
$wsOverview.Range("F6").Value = "Ready for handoff"
$wsOverview.Range("G6").Value = "2016-10-20 08:11:09"

$wsOverview.Range("A7").Value = "a35e0bd0-3195-4ccf-b9d6-de4f2dd34f28.md"
$wsOverview.Range("B7").Value = "e2e\a35e0bd0-3195-4ccf-b9d6-de4f2dd34f28.md"
$wsOverview.Range("E7").Value = "Ready for handoff"
$wsOverview.Range("F7").Value = "Ready for handoff"
$wsOverview.Range("G7").Value = "2016-10-20 08:17:51"

# Rebuild the hyperlinks on the "Overview" sheet so B6/B7 show the swapped
# display text while keeping every target URL as it was.
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/55e79d23f0cdc095bd1184b0f1770a5b49d12444/e2e/31a93d9a-6c2d-4f95-9781-35c4823c83b0.md", "", "", "e2e\31a93d9a-6c2d-4f95-9781-35c4823c83b0.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1ec3096b9d04b8c70211aa28e7c3772e1701bd18/e2e/417653bc-5fed-4be2-88eb-b0332b5678c9.md", "", "", "e2e\417653bc-5fed-4be2-88eb-b0332b5678c9.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1ec3096b9d04b8c70211aa28e7c3772e1701bd18/e2e/487d037c-ffb2-4963-a725-8a375519f0e8.md", "", "", "e2e\487d037c-ffb2-4963-a725-8a375519f0e8.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f40296a0935cfa95d017a82802fe941bcd3405c8/e2e/5811138a-fefe-4aa0-bb66-d88c61e7e508.md", "", "", "e2e\5811138a-fefe-4aa0-bb66-d88c61e7e508.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7a4f4df01daf117065603a0f69978dba50fd36b2/e2e/4d8f17aa-978e-4fa6-a9bf-edb8ba655f6c.md", "", "", "e2e\4d8f17aa-978e-4fa6-a9bf-edb8ba655f6c.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/70f2e4e716288b4cd75771faf375696e849d3899/e2e/a35e0bd0-3195-4ccf-b9d6-de4f2dd34f28.md", "", "", "e2e\a35e0bd0-3195-4ccf-b9d6-de4f2dd34f28.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn": columns A..P
#   A Source File Name, C Status, G Latest Handoff File, H Latest Handoff Datetime, I5 ...
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A6").Value = "4d8f17aa-978e-4fa6-a9bf-edb8ba655f6c.md"
$wsZh.Range("C6").Value = "Ready for handoff"
$wsZh.Range("G6").Value = "4d8f17aa-978e-4fa6-a9bf-edb8ba655f6c.f3621c206cbc6d5cd2d74b72e86bdd88c4a4c3bd.zh-cn.xlf"
$wsZh.Range("H6").Value = "2016-10-20 08:10:57"

$wsZh.Range("A7").Value = "a35e0bd0-3195-4ccf-b9d6-de4f2dd34f28.md"
$wsZh.Range("C7").Value = "Ready for handoff"
$wsZh.Range("G7").Value = "a35e0bd0-3195-4ccf-b9d6-de4f2dd34f28.d79cfe6515e087a2326730baea48cf32efd85fd9.zh-cn.xlf"
$wsZh.Range("H7").Value = "2016-10-20 08:17:38"

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/55e79d23f0cdc095bd1184b0f1770a5b49d12444/e2e/31a93d9a-6c2d-4f95-9781-35c4823c83b0.md", "", "", "31a93d9a-6c2d-4f95-9781-35c4823c83b0.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1ec3096b9d04b8c70211aa28e7c3772e1701bd18/e2e/417653bc-5fed-4be2-88eb-b0332b5678c9.md", "", "", "417653bc-5fed-4be2-88eb-b0332b5678c9.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1ec3096b9d04b8c70211aa28e7c3772e1701bd18/e2e/487d037c-ffb2-4963-a725-8a375519f0e8.md", "", "", "487d037c-ffb2-4963-a725-8a375519f0e8.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f40296a0935cfa95d017a82802fe941bcd3405c8/e2e/5811138a-fefe-4aa0-bb66-d88c61e7e508.md", "", "", "5811138a-fefe-4aa0-bb66-d88c61e7e508.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I5"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/ec340342b6df4f242c627e455556350410ca6364/e2e/5811138a-fefe-4aa0-bb66-d88c61e7e508.md", "", "", "5811138a-fefe-4aa0-bb66-d88c61e7e508.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7a4f4df01daf117065603a0f69978dba50fd36b2/e2e/4d8f17aa-978e-4fa6-a9bf-edb8ba655f6c.md", "", "", "4d8f17aa-978e-4fa6-a9bf-edb8ba655f6c.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/70f2e4e716288b4cd75771faf375696e849d3899/e2e/a35e0bd0-3195-4ccf-b9d6-de4f2dd34f28.md", "", "", "a35e0bd0-3195-4ccf-b9d6-de4f2dd34f28.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de": columns A..P (same layout as zh-cn)
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A6").Value = "4d8f17aa-978e-4fa6-a9bf-edb8ba655f6c.md"
$wsDe.Range("C6").Value = "Ready for handoff"
$wsDe.Range("G6").Value = "4d8f17aa-978e-4fa6-a9bf-edb8ba655f6c.f3621c206cbc6d5cd2d74b72e86bdd88c4a4c3bd.de-de.xlf"
$wsDe.Range("H6").Value = "2016-10-20 08:11:09"

$wsDe.Range("A7").Value = "a35e0bd0-3195-4ccf-b9d6-de4f2dd34f28.md"
$wsDe.Range("C7").Value = "Ready for handoff"
$wsDe.Range("G7").Value = "a35e0bd0-3195-4ccf-b9d6-de4f2dd34f28.d79cfe6515e087a2326730baea48cf32efd85fd9.de-de.xlf"
$wsDe.Range("H7").Value = "2016-10-20 08:17:51"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/55e79d23f0cdc095bd1184b0f1770a5b49d12444/e2e/31a93d9a-6c2d-4f95-9781-35c4823c83b0.md", "", "", "31a93d9a-6c2d-4f95-9781-35c4823c83b0.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1ec3096b9d04b8c70211aa28e7c3772e1701bd18/e2e/417653bc-5fed-4be2-88eb-b0332b5678c9.md", "", "", "417653bc-5fed-4be2-88eb-b0332b5678c9.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1ec3096b9d04b8c70211aa28e7c3772e1701bd18/e2e/487d037c-ffb2-4963-a725-8a375519f0e8.md", "", "", "487d037c-ffb2-4963-a725-8a375519f0e8.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f40296a0935cfa95d017a82802fe941bcd3405c8/e2e/5811138a-fefe-4aa0-bb66-d88c61e7e508.md", "", "", "5811138a-fefe-4aa0-bb66-d88c61e7e508.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I5"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/c97a49c374990eb845e981497e181c0e08cc3761/e2e/5811138a-fefe-4aa0-bb66-d88c61e7e508.md", "", "", "5811138a-fefe-4aa0-bb66-d88c61e7e508.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7a4f4df01daf117065603a0f69978dba50fd36b2/e2e/4d8f17aa-978e-4fa6-a9bf-edb8ba655f6c.md", "", "", "4d8f17aa-978e-4fa6-a9bf-edb8ba655f6c.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/70f2e4e716288b4cd75771faf375696e849d3899/e2e/a35e0bd0-3195-4ccf-b9d6-de4f2dd34f28.md", "", "", "a35e0bd0-3195-4ccf-b9d6-de4f2dd34f28.md") | Out-Null
